$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.501.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.203.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  -5.11%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.390"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.762.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "65.423.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.226.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "413.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.734.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.712"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0637"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "297.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.908"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
